$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Results for TestCase_F6 / TestCase_F7 moved from PASS to SKIP
$ws.Range("E7").Value = "SKIP"
$ws.Range("E8").Value = "SKIP"

# Drop the stray "apply fill" flag left on the Runmode cells for those two
# rows so their formatting matches the rest of the Runmode column.
$ws.Range("D7").Interior.Pattern = -4142
$ws.Range("D8").Interior.Pattern = -4142

# Clone the bordered row formatting down onto the two new rows before
# filling in their values.
$ws.Range("A6:E6").Copy($ws.Range("A9:E9"))
$ws.Range("A6:E6").Copy($ws.Range("A10:E10"))

# New test case: TestCase_F8
$ws.Range("A9").Value = "TestCase_F8"
$ws.Range("B9").Value = "OPQA-215"
$ws.Range("C9").Value = "Verify that user able to recevies a notification when other user commented on his post"
$ws.Range("D9").Value = "Y"
$ws.Range("E9").Value = "SKIP"

# New test case: TestCase_F9
$ws.Range("A10").Value = "TestCase_F9"
$ws.Range("B10").Value = "OPQA-216"
$ws.Range("C10").Value = "Verify that user receives a notification when someone he is following user comments on a post"
$ws.Range("D10").Value = "Y"
$ws.Range("E10").Value = "PASS"
